$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$apos = [string][char]39

$ws.Range("D2").Value = '65.439.90'
$ws.Range("E2").Value = '  +4.91%  '
$ws.Range("D3").Value = '3.499.26'
$ws.Range("E3").Value = '  +1.87%  '
$ws.Range("D4").Formula = $apos + '0.999'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Formula = $apos + '417.30'
$ws.Range("E5").Value = '  +0.81%  '
$ws.Range("E6").Value = '  +1.37%  '
$ws.Range("D7").Formula = $apos + '0.655'
$ws.Range("E7").Value = '  +5.10%  '
$ws.Range("D8").Formula = $apos + '0.999'
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Formula = $apos + '0.782'
$ws.Range("E9").Value = '  +7.41%  '
$ws.Range("E10").Value = '  +16.38%  '
$ws.Range("D11").Formula = $apos + '43.33'
$ws.Range("E11").Value = '  +0.97%  '
$ws.Range("D12").Formula = $apos + '0.0000268'
$ws.Range("E12").Value = '  +22.96%  '
$ws.Range("D13").Formula = $apos + '9.99'
$ws.Range("E13").Value = '  +8.56%  '
$ws.Range("D14").Value = '4.047.66'
$ws.Range("E14").Value = '  +1.80%  '
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("D16").Formula = $apos + '20.43'
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("D17").Value = '3.508.78'
$ws.Range("E17").Value = '  +2.79%  '
$ws.Range("D18").Formula = $apos + '12.88'
$ws.Range("E18").Value = '  +1.79%  '
$ws.Range("E19").Value = '  +2.20%  '
$ws.Range("D20").Value = '65.148.15'
$ws.Range("E20").Value = '  +4.38%  '
$ws.Range("D21").Formula = $apos + '445.48'
$ws.Range("E21").Value = '  -6.08%  '
$ws.Range("D22").Formula = $apos + '89.60'
$ws.Range("E22").Value = '  -2.17%  '
$ws.Range("D23").Formula = $apos + '3.24'
$ws.Range("E23").Value = '  -0.48%  '
$ws.Range("D24").Formula = $apos + '13.14'
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("E25").Value = '  +1.14%  '
$ws.Range("D26").Formula = $apos + '9.89'
$ws.Range("E26").Value = '  -1.19%  '
$ws.Range("D27").Formula = $apos + '33.98'
$ws.Range("E27").Value = '  +2.13%  '
$ws.Range("D28").Formula = $apos + '12.46'
$ws.Range("E28").Value = '  +4.88%  '
$ws.Range("D29").Formula = $apos + '2.74'
$ws.Range("E29").Value = '  +4.28%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Formula = $apos + '0.117'
$ws.Range("E30").Value = '  +5.12%  '
$ws.Range("B31").Value = 'RenderToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D31").Formula = $apos + '7.38'
$ws.Range("E31").Value = '  -5.11%  '
$ws.Range("D32").Formula = $apos + '0.163'
$ws.Range("E32").Value = '  -1.80%  '
$ws.Range("D33").Formula = $apos + '0.998'
$ws.Range("E33").Value = '  -0.13%  '
$ws.Range("D34").Formula = $apos + '39.41'
$ws.Range("E34").Value = '  -4.14%  '
$ws.Range("D35").Formula = $apos + '57.43'
$ws.Range("E35").Value = '  -0.63%  '
$ws.Range("E36").Value = '  +3.26%  '
$ws.Range("D37").Value = '0.0₃0726'
$ws.Range("E37").Value = '  +35.94%  '
$ws.Range("D38").Formula = $apos + '0.146'
$ws.Range("E38").Value = '  +8.60%  '
$ws.Range("D39").Formula = $apos + '0.996'
$ws.Range("E39").Value = '  -0.27%  '
$ws.Range("D40").Formula = $apos + '2.80'
$ws.Range("E40").Value = '  +5.62%  '
$ws.Range("E41").Value = '  -1.13%  '
$ws.Range("D42").Formula = $apos + '4.50'
$ws.Range("E42").Value = '  +3.83%  '
$ws.Range("D43").Formula = $apos + '147.28'
$ws.Range("E43").Value = '  +2.36%  '
$ws.Range("D44").Formula = $apos + '3.25'
$ws.Range("E44").Value = '  -2.39%  '
$ws.Range("E45").Value = '  -5.08%  '
$ws.Range("D46").Formula = $apos + '2.00'
$ws.Range("E46").Value = '  -3.62%  '
$ws.Range("D47").Formula = $apos + '2.32'
$ws.Range("E47").Value = '  -5.03%  '
$ws.Range("E48").Value = '  +5.25%  '
$ws.Range("D49").Formula = $apos + '15.72'
$ws.Range("E49").Value = '  -4.44%  '
$ws.Range("E50").Value = '  +8.23%  '
$ws.Range("D51").Formula = $apos + '21.51'
$ws.Range("E51").Value = '  -3.31%  '